# Add a new weekly record at the top of the price history block (row 534),
# pushing all existing rows 534:649 down to 535:650.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("534:534").Insert()

$ws.Range("A534").Value = 8
$ws.Range("B534").Value = "Terminal La Palmera de La Serena"
$ws.Range("C534").Value = "Coquimbo"
$ws.Range("D534").Value = 45275
$ws.Range("E534").Value = 4
$ws.Range("F534").Value = 100114013
$ws.Range("G534").Value = "Zanahoria"
$ws.Range("H534").Value = "Sin especificar"
$ws.Range("I534").Value = "Primera"
$ws.Range("J534").Value = 480
$ws.Range("K534").Value = 5500
$ws.Range("L534").Value = 6000
$ws.Range("M534").Value = 5750
$ws.Range("N534").Value = "`$/saco 20 kilos"
$ws.Range("O534").Value = "Provincia del Elquí"
$ws.Range("P534").Value = 288
$ws.Range("Q534").Value = 20
$ws.Range("R534").Value = "Hortaliza"
